$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the skill ratings for row 10 (NAYEE MUDDIN KHAN / MSCV / Project Management)
$ws.Range("D10").Value = "3 - Good"
$ws.Range("E10").Value = "3 - Good"
$ws.Range("F10").Value = "1 - Beginner"
$ws.Range("G10").Value = "3 - Good"
$ws.Range("H10").Value = "1 - Beginner"
$ws.Range("I10").Value = "3 - Good"

# Update the active selection to G7
$ws.Range("G7").Select()
